$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.144.68"
$ws.Range("E2").Value = "  +1.37%  "
$ws.Range("D3").Value = "2.589.31"
$ws.Range("E3").Value = "  +0.15%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'562.91"
$ws.Range("E5").Value = "  -0.44%  "
$ws.Range("D6").Value = "'141.62"
$ws.Range("E6").Value = "  -0.97%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "'0.596"
$ws.Range("E8").Value = "  -0.57%  "
$ws.Range("D9").Value = "2.612.60"
$ws.Range("E9").Value = "  +0.81%  "
$ws.Range("D10").Value = "'6.54"
$ws.Range("E10").Value = "  -1.64%  "
$ws.Range("E11").Value = "  +0.83%  "
$ws.Range("D12").Value = "'0.368"
$ws.Range("E12").Value = "  +6.50%  "
$ws.Range("E13").Value = "  -4.88%  "
$ws.Range("D14").Value = "3.055.88"
$ws.Range("E14").Value = "  +0.58%  "
$ws.Range("D15").Value = "60.184.94"
$ws.Range("E15").Value = "  +1.51%  "
$ws.Range("D16").Value = "'23.22"
$ws.Range("E16").Value = "  +2.45%  "
$ws.Range("D17").Value = "'0.0000139"
$ws.Range("E17").Value = "  +1.78%  "
$ws.Range("D18").Value = "2.606.10"
$ws.Range("E18").Value = "  +1.00%  "
$ws.Range("D19").Value = "'11.22"
$ws.Range("E19").Value = "  +8.05%  "
$ws.Range("D20").Value = "'4.63"
$ws.Range("E20").Value = "  +1.46%  "
$ws.Range("D21").Value = "'345.23"
$ws.Range("E21").Value = "  +2.26%  "
$ws.Range("D22").Value = "'6.95"
$ws.Range("E22").Value = "  +8.81%  "
$ws.Range("D23").Value = "'0.998"
$ws.Range("E23").Value = "  -0.25%  "
$ws.Range("D24").Value = "'0.530"
$ws.Range("E24").Value = "  +15.95%  "
$ws.Range("D25").Value = "'63.08"
$ws.Range("E25").Value = "  -1.92%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("E27").Value = "  -2.07%  "
$ws.Range("D28").Value = "'7.61"
$ws.Range("E28").Value = "  +3.82%  "
$ws.Range("D29").Value = "0.0₃0781"
$ws.Range("E29").Value = "  +0.31%  "
$ws.Range("D30").Value = "'1.78"
$ws.Range("E30").Value = "  +6.40%  "
$ws.Range("E31").Value = "  -0.09%  "
$ws.Range("D32").Value = "'6.30"
$ws.Range("E32").Value = "  +3.59%  "
$ws.Range("D33").Value = "'160.92"
$ws.Range("E33").Value = "  -0.55%  "
$ws.Range("D34").Value = "'19.39"
$ws.Range("E34").Value = "  +2.35%  "
$ws.Range("D35").Value = "'4.20"
$ws.Range("E35").Value = "  +4.52%  "
$ws.Range("D36").Value = "'0.952"
$ws.Range("E36").Value = "  +8.37%  "
$ws.Range("D37").Value = "'1.21"
$ws.Range("E37").Value = "  +3.69%  "
$ws.Range("D38").Value = "'1.59"
$ws.Range("E38").Value = "  +7.09%  "
$ws.Range("D39").Value = "'37.68"
$ws.Range("E39").Value = "  +0.33%  "
$ws.Range("D40").Value = "'0.854"
$ws.Range("E40").Value = "  -2.88%  "
$ws.Range("D41").Value = "'3.80"
$ws.Range("E41").Value = "  +3.74%  "
$ws.Range("D42").Value = "'292.81"
$ws.Range("E42").Value = "  -2.25%  "
$ws.Range("D43").Value = "'137.24"
$ws.Range("E43").Value = "  +3.97%  "
$ws.Range("D44").Value = "'0.999"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").Value = "'0.604"
$ws.Range("E45").Value = "  +1.16%  "
$ws.Range("D46").Value = "'0.0976"
$ws.Range("E46").Value = "  +0.21%  "
$ws.Range("D47").Value = "'0.0545"
$ws.Range("E47").Value = "  +1.64%  "
$ws.Range("D48").Value = "'19.44"
$ws.Range("E48").Value = "  +1.52%  "
$ws.Range("E49").Value = "  +2.90%  "
$ws.Range("D50").Value = "'10.68"
$ws.Range("E50").Value = "  +0.53%  "
$ws.Range("D51").Value = "'19.57"
$ws.Range("E51").Value = "  +5.67%  "
